$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 1136
$ws1.Range("F11").Value = 782
$ws1.Range("F12").Value = 68

# Sheet "全部类型": same rows duplicated, update accordingly
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 1136
$ws4.Range("F12").Value = 782
$ws4.Range("F13").Value = 68
